# Update "想去人数" (F) and "最低票价" (G) figures on the "展览" and "全部类型" sheets
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5407
$ws1.Range("F3").Value = 589
$ws1.Range("F4").Value = 11574
$ws1.Range("G4").Value = 58
$ws1.Range("F5").Value = 284
$ws1.Range("F6").Value = 592
$ws1.Range("F8").Value = 263
$ws1.Range("F9").Value = 1005
$ws1.Range("F10").Value = 97

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 5407
$ws4.Range("F5").Value = 589
$ws4.Range("F7").Value = 11574
$ws4.Range("G7").Value = 58
$ws4.Range("F8").Value = 284
$ws4.Range("F9").Value = 592
$ws4.Range("F13").Value = 263
$ws4.Range("F14").Value = 1005
$ws4.Range("F16").Value = 97
